$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate (by content, before any mutation) the two paragraphs near the end of
# the document that this edit touches:
#   - the bold "Play Almighty Sparta Free..." paragraph (duplicate title),
#     which gets removed from the bottom and re-purposed as the new
#     "Meta description" paragraph right under the real title, and
#   - the italic paragraph, whose text is replaced with the DALLE prompt.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldIdx = -1
$italicIdx = -1
for ($i = $count; $i -ge 2; $i--) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($boldIdx -eq -1 -and $t -eq "Play Almighty Sparta Free: Game Review and Pros & Cons") {
        $boldIdx = $i
    }
    if ($italicIdx -eq -1 -and $t -eq "Read our review of Almighty Sparta online slot game. Learn about its pros & cons before you play for free. Discover if the game is worth your time.") {
        $italicIdx = $i
    }
}

# ---------------------------------------------------------------------------
# 1) Insert a new empty paragraph right after the H1 title paragraph.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(2)
$newPara.Style = "Normal"

# ---------------------------------------------------------------------------
# 2) Clone the run layout (leading empty run + bold run) of the bold
#    "Play Almighty Sparta Free..." paragraph into the new paragraph, by
#    copying its FormattedText, then retarget the bold run's text to
#    "Meta description" and append the rest of the sentence as a plain run.
# ---------------------------------------------------------------------------
$boldIdx = $boldIdx + 1   # shifted down by one because of the insert above
$srcPara = $d.Paragraphs($boldIdx)
$newPara.Range.FormattedText = $srcPara.Range.FormattedText

$newPara = $d.Paragraphs(2)
$boldRunRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$boldRunRange.Text = "Meta description"

$newPara = $d.Paragraphs(2)
$newPara.Range.InsertAfter(": Read our review of Almighty Sparta online slot game. Learn about its pros & cons before you play for free. Discover if the game is worth your time.")

# ---------------------------------------------------------------------------
# 3) Remove the original bold "Play Almighty Sparta Free..." paragraph.
# ---------------------------------------------------------------------------
$d.Paragraphs($boldIdx).Range.Delete()

# ---------------------------------------------------------------------------
# 4) Replace the italic paragraph's text with the DALLE prompt, keeping its
#    existing run formatting (italic).
# ---------------------------------------------------------------------------
$oldText = "Read our review of Almighty Sparta online slot game. Learn about its pros & cons before you play for free. Discover if the game is worth your time."
$newText = "Prompt for DALLE: Create a vibrant feature image for Almighty Sparta slot game that catches the eye of online casino players. The image should be in cartoon style, featuring a happy Maya warrior with glasses. The Maya warrior should be standing in an epic pose, wearing a helmet, a shield on one arm, and holding a sword in the other hand. The background should include a scenic view of ancient Sparta with the sea, mountains, and clouds. Use vibrant colors such as golden yellow, olive green, and sky blue to depict the glorious historical era of ancient Sparta. Make sure the image has an appealing design and an eye-catching layout that represents the adventurous and thrilling gameplay of the slot game."

$italicRange = $d.Paragraphs($italicIdx).Range
$italicRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
